# Fruta / hortaliza, semanal
# Insert a new week's worth of price data (4 quality-grade rows, date 2021-11-09 /
# serial 44509) at the top of the "Pina - Caramelo" block on Sheet1, pushing the
# existing historical rows down by 4 (Excel copies formatting/styles automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 770; rows 770-797 shift down to 774-801.
$ws.Rows("770:773").Insert()

# Row 770: Especial
$ws.Range("A770").Value = 6
$ws.Range("B770").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C770").Value = "Metropolitana"
$ws.Range("D770").Value = 44509
$ws.Range("E770").Value = 13
$ws.Range("F770").Value = "Fruta"
$ws.Range("G770").Value = 100108
$ws.Range("H770").Value = "Tropicales y subtropicales"
$ws.Range("I770").Value = 100108005
$ws.Range("J770").Value = "Piña"
$ws.Range("K770").Value = "Caramelo"
$ws.Range("L770").Value = "Especial"
$ws.Range("M770").Value = 368
$ws.Range("N770").Value = 17000
$ws.Range("O770").Value = 18000
$ws.Range("P770").Value = 17853
$ws.Range("Q770").Value = "$/caja 10 unidades"
$ws.Range("R770").Value = "Ecuador"
$ws.Range("S770").Value = 1785
$ws.Range("T770").Value = 10

# Row 771: Primera
$ws.Range("A771").Value = 6
$ws.Range("B771").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C771").Value = "Metropolitana"
$ws.Range("D771").Value = 44509
$ws.Range("E771").Value = 13
$ws.Range("F771").Value = "Fruta"
$ws.Range("G771").Value = 100108
$ws.Range("H771").Value = "Tropicales y subtropicales"
$ws.Range("I771").Value = 100108005
$ws.Range("J771").Value = "Piña"
$ws.Range("K771").Value = "Caramelo"
$ws.Range("L771").Value = "Primera"
$ws.Range("M771").Value = 506
$ws.Range("N771").Value = 17000
$ws.Range("O771").Value = 18000
$ws.Range("P771").Value = 17787
$ws.Range("Q771").Value = "$/caja 12 unidades"
$ws.Range("R771").Value = "Ecuador"
$ws.Range("S771").Value = 1482
$ws.Range("T771").Value = 12

# Row 772: Segunda
$ws.Range("A772").Value = 6
$ws.Range("B772").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C772").Value = "Metropolitana"
$ws.Range("D772").Value = 44509
$ws.Range("E772").Value = 13
$ws.Range("F772").Value = "Fruta"
$ws.Range("G772").Value = 100108
$ws.Range("H772").Value = "Tropicales y subtropicales"
$ws.Range("I772").Value = 100108005
$ws.Range("J772").Value = "Piña"
$ws.Range("K772").Value = "Caramelo"
$ws.Range("L772").Value = "Segunda"
$ws.Range("M772").Value = 486
$ws.Range("N772").Value = 17000
$ws.Range("O772").Value = 18000
$ws.Range("P772").Value = 17778
$ws.Range("Q772").Value = "$/caja 14 unidades"
$ws.Range("R772").Value = "Ecuador"
$ws.Range("S772").Value = 1270
$ws.Range("T772").Value = 14

# Row 773: Tercera
$ws.Range("A773").Value = 6
$ws.Range("B773").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C773").Value = "Metropolitana"
$ws.Range("D773").Value = 44509
$ws.Range("E773").Value = 13
$ws.Range("F773").Value = "Fruta"
$ws.Range("G773").Value = 100108
$ws.Range("H773").Value = "Tropicales y subtropicales"
$ws.Range("I773").Value = 100108005
$ws.Range("J773").Value = "Piña"
$ws.Range("K773").Value = "Caramelo"
$ws.Range("L773").Value = "Tercera"
$ws.Range("M773").Value = 281
$ws.Range("N773").Value = 17000
$ws.Range("O773").Value = 18000
$ws.Range("P773").Value = 17616
$ws.Range("Q773").Value = "$/caja 16 unidades"
$ws.Range("R773").Value = "Ecuador"
$ws.Range("S773").Value = 1101
$ws.Range("T773").Value = 16
